$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    # Force the cell to keep its string content as text (matching the
    # original inlineStr cell type) instead of Excel auto-coercing
    # numeric-looking strings (e.g. "1.001") into numbers, then restore
    # the default "Normal" style so no stray formatting is introduced.
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "28.115.59"
Set-TextCell "E2" "  -1.64%  "

Set-TextCell "D3" "1.800.24"
Set-TextCell "E3" "  +0.14%  "

Set-TextCell "D4" "1.001"
Set-TextCell "E4" "  -0.01%  "

Set-TextCell "D5" "316.04"
Set-TextCell "E5" "  +0.85%  "

Set-TextCell "D6" "1.001"
Set-TextCell "E6" "  +0.10%  "

Set-TextCell "D7" "0.5452"
Set-TextCell "E7" "  +1.46%  "

Set-TextCell "D8" "0.3783"
Set-TextCell "E8" "  -0.04%  "

Set-TextCell "D9" "0.07457"
Set-TextCell "E9" "  -1.12%  "

Set-TextCell "D10" "41.98"
Set-TextCell "E10" "  -1.37%  "

Set-TextCell "E11" "  -2.09%  "

Set-TextCell "D12" "1.001"
Set-TextCell "E12" "  +0.20%  "

Set-TextCell "D13" "6.200"
Set-TextCell "E13" "  +0.17%  "

Set-TextCell "D14" "20.43"
Set-TextCell "E14" "  -2.76%  "

Set-TextCell "D15" "7.337"
Set-TextCell "E15" "  -1.62%  "

Set-TextCell "D16" "1.798.20"
Set-TextCell "E16" "  +0.06%  "

Set-TextCell "D17" "89.63"
Set-TextCell "E17" "  -0.78%  "

Set-TextCell "D18" "0.00001063"
Set-TextCell "E18" "  -0.25%  "

Set-TextCell "D19" "0.06539"
Set-TextCell "E19" "  +1.34%  "

Set-TextCell "D20" "17.44"
Set-TextCell "E20" "  +1.24%  "

Set-TextCell "E21" "  +0.14%  "

Set-TextCell "D22" "5.930"
Set-TextCell "E22" "  +0.07%  "

Set-TextCell "D23" "28.142.98"
Set-TextCell "E23" "  -1.68%  "

Set-TextCell "E24" "  +0.10%  "

Set-TextCell "D25" "2.089"
Set-TextCell "E25" "  -0.38%  "

Set-TextCell "D26" "155.44"
Set-TextCell "E26" "  -3.37%  "

Set-TextCell "D27" "20.41"
Set-TextCell "E27" "  -0.24%  "

Set-TextCell "D28" "2.006.05"
Set-TextCell "E28" "  +0.10%  "

Set-TextCell "D29" "2.324"
Set-TextCell "E29" "  -2.27%  "

Set-TextCell "D30" "121.69"
Set-TextCell "E30" "  -1.26%  "

Set-TextCell "D31" "0.1112"
Set-TextCell "E31" "  +8.28%  "

Set-TextCell "D32" "1.116"
Set-TextCell "E32" "  +0.66%  "

Set-TextCell "D33" "3.665"
Set-TextCell "E33" "  -0.48%  "

Set-TextCell "D34" "5.559"
Set-TextCell "E34" "  -2.10%  "

Set-TextCell "D35" "0.06948"
Set-TextCell "E35" "  +7.28%  "

Set-TextCell "D36" "0.2221"
Set-TextCell "E36" "  -1.80%  "

Set-TextCell "D37" "0.02295"
Set-TextCell "E37" "  -0.40%  "

Set-TextCell "D38" "5.077"
Set-TextCell "E38" "  +0.46%  "

Set-TextCell "D39" "8.455"
Set-TextCell "E39" "  -4.98%  "

Set-TextCell "E40" "  -1.76%  "

Set-TextCell "D41" "0.6158"
Set-TextCell "E41" "  -1.65%  "

Set-TextCell "E42" "  +2.47%  "

Set-TextCell "E43" "  -3.26%  "

Set-TextCell "D44" "13.31"
Set-TextCell "E44" "  -1.15%  "

Set-TextCell "E45" "  +0.54%  "

Set-TextCell "D46" "0.5736"
Set-TextCell "E46" "  -2.69%  "

Set-TextCell "D47" "124.51"
Set-TextCell "E47" "  -1.01%  "

Set-TextCell "D48" "1.183"
Set-TextCell "E48" "  +2.04%  "

Set-TextCell "D49" "1.917"
Set-TextCell "E49" "  -2.64%  "

Set-TextCell "D50" "0.06816"
Set-TextCell "E50" "  -1.53%  "

Set-TextCell "D51" "71.79"
Set-TextCell "E51" "  -1.27%  "

